$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Change 1: move the hidden "_GoBack" bookmark from the end of the
# "SOURCES += main.cpp" paragraph to the end of the
# "TARGET = main  # Hier wird der Zielname festgelegt" paragraph. ---

$pSources = $d.Paragraphs.Item(5)
$xmlSources = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:firstLine="708"/></w:pPr>' + `
  '<w:r><w:t>SOURCES += main.cpp</w:t></w:r></w:p>'
$pSources.Range.InsertXML($xmlSources)

$pTarget = $d.Paragraphs.Item(6)
$xmlTarget = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:ind w:firstLine="708"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">TARGET = </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>main</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">  #</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Hier wird der </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Zielname</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> festgelegt</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pTarget.Range.InsertXML($xmlTarget)

# --- Change 2: add a new bullet after the "...main.pro" paragraph
# describing the Qt Creator build alternative. ---

$pPro = $d.Paragraphs.Item(7)
$newPara = $pPro.Range
$newPara.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item(8)

$xmlNew = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Listenabsatz"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Oder &#252;ber den </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Qt</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Creator das Projekt &#246;ffnen und </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Debug</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">/Release </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>build</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> erstellen.</w:t></w:r></w:p>'
$pNew.Range.InsertXML($xmlNew)

# Give the whole new paragraph's text the "Fett" (Strong) character style,
# matching every run in the source paragraph.
$pNew2 = $d.Paragraphs.Item(8)
$textStart = $pNew2.Range.Start
$textEnd = $pNew2.Range.End - 1
$fullTextRange = $d.Range($textStart, $textEnd)
$fullTextRange.Style = "Fett"
